# Generate report for ARM
# Update the ARM row (row 3) metrics in columns H:Q, then refresh the
# sheet's scroll/selection state to match where the user ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- ARM row (row 3) value updates ---------------------------------------
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 5.5
$ws.Range("J3:M3").ClearContents()
$ws.Range("N3").Value = 5.8
$ws.Range("O3").Value = 6.2
$ws.Range("P3").Value = 1.9
$ws.Range("Q3").Value = 3.4

# --- Sheet view: scroll frozen pane back to top, move the selection -------
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("N4").Select() | Out-Null
